$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.251.42'
$ws.Range("E2").Value = '  -0.09%  '
$ws.Range("D3").Value = '2.005.16'
$ws.Range("E3").Value = '  -1.23%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '258.69'
$ws.Range("E5").Value = '  +4.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.617'
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '56.52'
$ws.Range("E8").Value = '  -6.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.381'
$ws.Range("E9").Value = '  -3.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0771'
$ws.Range("E10").Value = '  -4.76%  '
$ws.Range("E11").Value = '  -3.20%  '
$ws.Range("D12").Value = '2.302.85'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.22'
$ws.Range("E13").Value = '  -6.88%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.68'
$ws.Range("E14").Value = '  -2.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.792'
$ws.Range("E15").Value = '  -7.70%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.20'
$ws.Range("E16").Value = '  -5.77%  '
$ws.Range("D17").Value = '2.010.23'
$ws.Range("E17").Value = '  -1.10%  '
$ws.Range("D18").Value = '37.295.72'
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.10'
$ws.Range("E19").Value = '  -0.77%  '
$ws.Range("D20").Value = '0.0₃0833'
$ws.Range("E20").Value = '  -3.68%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '233.54'
$ws.Range("E21").Value = '  +1.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.09'
$ws.Range("E22").Value = '  -3.03%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.57'
$ws.Range("E24").Value = '  -1.62%  '
$ws.Range("E25").Value = '  -0.42%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.79'
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.94'
$ws.Range("E27").Value = '  -5.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.52'
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E29").Value = '  -7.21%  '
$ws.Range("E30").Value = '  -4.11%  '
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.58'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0638'
$ws.Range("E33").Value = '  -5.44%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.44'
$ws.Range("E34").Value = '  -1.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.35'
$ws.Range("E35").Value = '  -5.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.41'
$ws.Range("E36").Value = '  -6.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.82'
$ws.Range("E37").Value = '  +0.44%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("E40").Value = '  +0.76%  '
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("E42").Value = '  -1.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0923'
$ws.Range("E43").Value = '  -5.85%  '
$ws.Range("D44").Value = '1.432.77'
$ws.Range("E44").Value = '  +3.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '89.13'
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '15.56'
$ws.Range("E46").Value = '  -9.27%  '
$ws.Range("E47").Value = '  -3.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.93'
$ws.Range("E48").Value = '  +2.61%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.96'
$ws.Range("E49").Value = '  -6.81%  '
$ws.Range("D50").Value = '2.193.77'
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.94'
$ws.Range("E51").Value = '  -10.01%  '
